# Update automàtic: dades i banners [2026-02-10 03:50]
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regular text/value updates (dates, measurements with units, etc.)
$ws.Cells.Item(2, 5).Value = "2026-02-10 03:48:28"
$ws.Cells.Item(2, 9).Value = "7.6 mm"
$ws.Cells.Item(3, 5).Value = "2026-02-10 03:48:30"
$ws.Cells.Item(3, 7).Value = "181 cm"
$ws.Cells.Item(3, 9).Value = "5.5 mm"
$ws.Cells.Item(4, 5).Value = "2026-02-10 03:48:33"
$ws.Cells.Item(4, 10).Value = "1004.9 hPa"
$ws.Cells.Item(5, 5).Value = "2026-02-10 03:48:35"
$ws.Cells.Item(5, 9).Value = "8.1 mm"
$ws.Cells.Item(5, 15).Value = "-0.4 °C"
$ws.Cells.Item(6, 5).Value = "2026-02-10 03:48:37"
$ws.Cells.Item(7, 5).Value = "2026-02-10 03:48:40"
$ws.Cells.Item(7, 10).Value = "1005.2 hPa"
$ws.Cells.Item(7, 13).Value = "12.4 °C 3:05 TU"
$ws.Cells.Item(7, 15).Value = "12.1 °C"
$ws.Cells.Item(8, 5).Value = "2026-02-10 03:48:42"
$ws.Cells.Item(8, 10).Value = "1005.3 hPa"
$ws.Cells.Item(9, 5).Value = "2026-02-10 03:48:44"
$ws.Cells.Item(9, 14).Value = "5.9 °C 3:25 TU"
$ws.Cells.Item(9, 15).Value = "7.0 °C"
$ws.Cells.Item(10, 5).Value = "2026-02-10 03:48:47"
$ws.Cells.Item(10, 14).Value = "6.3 °C 3:27 TU"
$ws.Cells.Item(10, 15).Value = "7.5 °C"
$ws.Cells.Item(11, 5).Value = "2026-02-10 03:48:49"
$ws.Cells.Item(11, 14).Value = "2.6 °C 3:00 TU"
$ws.Cells.Item(12, 5).Value = "2026-02-10 03:48:51"
$ws.Cells.Item(12, 14).Value = "6.5 °C 3:27 TU"
$ws.Cells.Item(12, 15).Value = "7.4 °C"
$ws.Cells.Item(13, 5).Value = "2026-02-10 03:48:53"
$ws.Cells.Item(13, 10).Value = "1008.0 hPa"
$ws.Cells.Item(13, 14).Value = "2.4 °C 3:16 TU"
$ws.Cells.Item(13, 15).Value = "2.7 °C"
$ws.Cells.Item(14, 5).Value = "2026-02-10 03:48:56"
$ws.Cells.Item(14, 14).Value = "8.7 °C 3:29 TU"
$ws.Cells.Item(14, 15).Value = "9.9 °C"
$ws.Cells.Item(15, 5).Value = "2026-02-10 03:48:58"
$ws.Cells.Item(15, 14).Value = "5.7 °C 3:28 TU"
$ws.Cells.Item(15, 15).Value = "7.4 °C"
$ws.Cells.Item(16, 5).Value = "2026-02-10 03:49:00"
$ws.Cells.Item(16, 9).Value = "5.8 mm"
$ws.Cells.Item(17, 5).Value = "2026-02-10 03:49:03"
$ws.Cells.Item(17, 13).Value = "2.6 °C 3:27 TU"
$ws.Cells.Item(17, 15).Value = "1.6 °C"
$ws.Cells.Item(18, 5).Value = "2026-02-10 03:49:05"
$ws.Cells.Item(18, 10).Value = "1005.1 hPa"
$ws.Cells.Item(18, 14).Value = "7.2 °C 3:09 TU"
$ws.Cells.Item(18, 15).Value = "8.0 °C"
$ws.Cells.Item(19, 5).Value = "2026-02-10 03:49:08"
$ws.Cells.Item(19, 13).Value = "4.3 °C 3:26 TU"
$ws.Cells.Item(19, 15).Value = "4.0 °C"
$ws.Cells.Item(20, 5).Value = "2026-02-10 03:49:10"
$ws.Cells.Item(20, 9).Value = "1.8 mm"
$ws.Cells.Item(20, 12).Value = "43.6 km/h - 343º 3:17 TU"
$ws.Cells.Item(20, 13).Value = "-0.6 °C 3:21 TU"
$ws.Cells.Item(20, 15).Value = "-1.5 °C"
$ws.Cells.Item(21, 5).Value = "2026-02-10 03:49:12"
$ws.Cells.Item(21, 14).Value = "4.1 °C 3:26 TU"
$ws.Cells.Item(22, 5).Value = "2026-02-10 03:49:15"
$ws.Cells.Item(22, 15).Value = "-2.1 °C"
$ws.Cells.Item(23, 5).Value = "2026-02-10 03:49:17"
$ws.Cells.Item(23, 7).Value = "180 cm"
$ws.Cells.Item(23, 9).Value = "5.9 mm"
$ws.Cells.Item(23, 15).Value = "-0.6 °C"
$ws.Cells.Item(24, 5).Value = "2026-02-10 03:49:20"
$ws.Cells.Item(24, 10).Value = "1006.8 hPa"
$ws.Cells.Item(24, 15).Value = "8.4 °C"
$ws.Cells.Item(25, 5).Value = "2026-02-10 03:49:22"
$ws.Cells.Item(25, 9).Value = "2.7 mm"
$ws.Cells.Item(25, 15).Value = "-0.9 °C"
$ws.Cells.Item(26, 5).Value = "2026-02-10 03:49:24"
$ws.Cells.Item(26, 14).Value = "2.2 °C 3:18 TU"
$ws.Cells.Item(26, 15).Value = "2.9 °C"
$ws.Cells.Item(27, 5).Value = "2026-02-10 03:49:27"
$ws.Cells.Item(27, 13).Value = "0.0 °C 3:18 TU"
$ws.Cells.Item(27, 15).Value = "-0.7 °C"
$ws.Cells.Item(28, 5).Value = "2026-02-10 03:49:29"
$ws.Cells.Item(28, 10).Value = "1005.4 hPa"
$ws.Cells.Item(28, 14).Value = "4.5 °C 3:27 TU"
$ws.Cells.Item(28, 15).Value = "5.6 °C"
$ws.Cells.Item(29, 5).Value = "2026-02-10 03:49:32"
$ws.Cells.Item(30, 5).Value = "2026-02-10 03:49:34"
$ws.Cells.Item(31, 5).Value = "2026-02-10 03:49:37"
$ws.Cells.Item(31, 14).Value = "8.8 °C 3:29 TU"
$ws.Cells.Item(32, 5).Value = "2026-02-10 03:49:39"
$ws.Cells.Item(32, 12).Value = "24.1 km/h - 304º 3:21 TU"
$ws.Cells.Item(32, 15).Value = "7.4 °C"
$ws.Cells.Item(33, 5).Value = "2026-02-10 03:49:41"
$ws.Cells.Item(33, 9).Value = "0.9 mm"
$ws.Cells.Item(33, 14).Value = "1.8 °C 3:12 TU"
$ws.Cells.Item(34, 5).Value = "2026-02-10 03:49:44"
$ws.Cells.Item(34, 9).Value = "1.3 mm"
$ws.Cells.Item(35, 5).Value = "2026-02-10 03:49:46"
$ws.Cells.Item(36, 5).Value = "2026-02-10 03:49:49"
$ws.Cells.Item(36, 12).Value = "39.2 km/h - 18º 3:08 TU"
$ws.Cells.Item(36, 13).Value = "10.8 °C 3:14 TU"
$ws.Cells.Item(36, 15).Value = "9.1 °C"
$ws.Cells.Item(37, 5).Value = "2026-02-10 03:49:51"
$ws.Cells.Item(37, 14).Value = "2.6 °C 3:24 TU"
$ws.Cells.Item(37, 15).Value = "4.3 °C"
$ws.Cells.Item(38, 5).Value = "2026-02-10 03:49:54"
$ws.Cells.Item(38, 12).Value = "21.2 km/h - 287º 3:29 TU"
$ws.Cells.Item(38, 13).Value = "8.6 °C 3:29 TU"
$ws.Cells.Item(38, 14).Value = "7.2 °C 3:19 TU"
$ws.Cells.Item(39, 5).Value = "2026-02-10 03:49:56"
$ws.Cells.Item(39, 9).Value = "1.0 mm"
$ws.Cells.Item(40, 5).Value = "2026-02-10 03:49:58"
$ws.Cells.Item(40, 9).Value = "1.6 mm"
$ws.Cells.Item(40, 15).Value = "4.9 °C"
$ws.Cells.Item(41, 5).Value = "2026-02-10 03:50:01"
$ws.Cells.Item(41, 10).Value = "1005.1 hPa"
$ws.Cells.Item(41, 14).Value = "8.9 °C 3:25 TU"
$ws.Cells.Item(41, 15).Value = "10.4 °C"
$ws.Cells.Item(42, 5).Value = "2026-02-10 03:50:03"
$ws.Cells.Item(43, 5).Value = "2026-02-10 03:50:05"
$ws.Cells.Item(43, 14).Value = "6.0 °C 3:13 TU"
$ws.Cells.Item(44, 5).Value = "2026-02-10 03:50:08"
$ws.Cells.Item(44, 7).Value = "217 cm"
$ws.Cells.Item(44, 9).Value = "5.7 mm"
$ws.Cells.Item(44, 15).Value = "-0.5 °C"
$ws.Cells.Item(45, 5).Value = "2026-02-10 03:50:10"
$ws.Cells.Item(45, 9).Value = "8.4 mm"
$ws.Cells.Item(46, 5).Value = "2026-02-10 03:50:12"
$ws.Cells.Item(46, 15).Value = "10.0 °C"

# Percentage-looking text values need to be forced to stay as text
# (Excel auto-converts "NN%" strings into numeric percentages otherwise).
# Setting NumberFormat to "@" (Text) first keeps the literal string.
$ws.Cells.Item(8, 8).NumberFormat = "@"
$ws.Cells.Item(8, 8).Value = "96%"
$ws.Cells.Item(18, 8).NumberFormat = "@"
$ws.Cells.Item(18, 8).Value = "97%"
$ws.Cells.Item(23, 8).NumberFormat = "@"
$ws.Cells.Item(23, 8).Value = "92%"
$ws.Cells.Item(25, 8).NumberFormat = "@"
$ws.Cells.Item(25, 8).Value = "92%"
$ws.Cells.Item(26, 8).NumberFormat = "@"
$ws.Cells.Item(26, 8).Value = "87%"
$ws.Cells.Item(31, 8).NumberFormat = "@"
$ws.Cells.Item(31, 8).Value = "85%"
$ws.Cells.Item(34, 8).NumberFormat = "@"
$ws.Cells.Item(34, 8).Value = "78%"
$ws.Cells.Item(36, 8).NumberFormat = "@"
$ws.Cells.Item(36, 8).Value = "94%"
$ws.Cells.Item(39, 8).NumberFormat = "@"
$ws.Cells.Item(39, 8).Value = "84%"
$ws.Cells.Item(40, 8).NumberFormat = "@"
$ws.Cells.Item(40, 8).Value = "97%"

# The NumberFormat="@" change above creates a new cell style, which would
# alter these cells' style index relative to the original workbook.
# Restore the original style (borders/alignment/General format) by copying
# the format from an untouched cell that already uses that same style,
# while keeping the text value we just set.
$styleRef = $ws.Cells.Item(2, 6)
$styleRef.Copy()
$ws.Cells.Item(8, 8).PasteSpecial(-4122)
$ws.Cells.Item(18, 8).PasteSpecial(-4122)
$ws.Cells.Item(23, 8).PasteSpecial(-4122)
$ws.Cells.Item(25, 8).PasteSpecial(-4122)
$ws.Cells.Item(26, 8).PasteSpecial(-4122)
$ws.Cells.Item(31, 8).PasteSpecial(-4122)
$ws.Cells.Item(34, 8).PasteSpecial(-4122)
$ws.Cells.Item(36, 8).PasteSpecial(-4122)
$ws.Cells.Item(39, 8).PasteSpecial(-4122)
$ws.Cells.Item(40, 8).PasteSpecial(-4122)
$excel.CutCopyMode = 0
